$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 44927
$ws.Cells.Item(4, 2).Value = 22403.435
$ws.Cells.Item(4, 3).Value = 300.356
$ws.Cells.Item(4, 4).Value = 3.4
$ws.Cells.Item(4, 5).Value = 4.33
$ws.Cells.Item(4, 6).Value = 3.53
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0.3838777983459662
$ws.Cells.Item(4, 12).Value = 301.203
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 15).Value = 5.882352941176472
$ws.Cells.Item(4, 16).Value = 3.5
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 19).Value = 5.542725173210172
$ws.Cells.Item(4, 20).Value = 4.516666666666667
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 23).Value = 6.232294617563738
$ws.Cells.Item(4, 24).Value = 3.646666666666667
$ws.Cells.Item(4, 25).Value = 0

$ws.Cells.Item(5, 1).Value = 44958
$ws.Cells.Item(5, 2).Value = 22403.435
$ws.Cells.Item(5, 3).Value = 301.509
$ws.Cells.Item(5, 4).Value = 3.6
$ws.Cells.Item(5, 5).Value = 4.57
$ws.Cells.Item(5, 6).Value = 3.75
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0.3838777983459662
$ws.Cells.Item(5, 12).Value = 301.203
$ws.Cells.Item(5, 13).Value = 0.3838777983459662
$ws.Cells.Item(5, 15).Value = 5.882352941176472
$ws.Cells.Item(5, 16).Value = 3.5
$ws.Cells.Item(5, 17).Value = 5.882352941176472
$ws.Cells.Item(5, 19).Value = 5.542725173210172
$ws.Cells.Item(5, 20).Value = 4.516666666666667
$ws.Cells.Item(5, 21).Value = 5.542725173210172
$ws.Cells.Item(5, 23).Value = 6.232294617563738
$ws.Cells.Item(5, 24).Value = 3.646666666666667
$ws.Cells.Item(5, 25).Value = 6.232294617563738

$ws.Cells.Item(6, 1).Value = 44986
$ws.Cells.Item(6, 2).Value = 22403.435
$ws.Cells.Item(6, 3).Value = 301.744
$ws.Cells.Item(6, 4).Value = 3.5
$ws.Cells.Item(6, 5).Value = 4.65
$ws.Cells.Item(6, 6).Value = 3.66
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 11).Value = 0.07794128865141747
$ws.Cells.Item(6, 12).Value = 301.203
$ws.Cells.Item(6, 13).Value = 0.4621182863002726
$ws.Cells.Item(6, 15).Value = -2.777777777777779
$ws.Cells.Item(6, 16).Value = 3.5
$ws.Cells.Item(6, 17).Value = 2.941176470588247
$ws.Cells.Item(6, 19).Value = 1.750547045951856
$ws.Cells.Item(6, 20).Value = 4.516666666666667
$ws.Cells.Item(6, 21).Value = 7.390300230946889
$ws.Cells.Item(6, 23).Value = -2.399999999999991
$ws.Cells.Item(6, 24).Value = 3.646666666666667
$ws.Cells.Item(6, 25).Value = 3.682719546742219

$ws.Cells.Item(7, 1).Value = 45017
$ws.Cells.Item(7, 2).Value = 22539.418
$ws.Cells.Item(7, 3).Value = 303.032
$ws.Cells.Item(7, 4).Value = 3.4
$ws.Cells.Item(7, 5).Value = 4.83
$ws.Cells.Item(7, 6).Value = 3.46
$ws.Cells.Item(7, 7).Value = 0.6069738859241891
$ws.Cells.Item(7, 9).Value = 0.6069738859241891
$ws.Cells.Item(7, 11).Value = 0.4268519009491323
$ws.Cells.Item(7, 12).Value = 302.095
$ws.Cells.Item(7, 13).Value = 0.8909427479391052
$ws.Cells.Item(7, 15).Value = -2.857142857142858
$ws.Cells.Item(7, 16).Value = 3.5
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(7, 19).Value = 3.870967741935472
$ws.Cells.Item(7, 20).Value = 4.683333333333334
$ws.Cells.Item(7, 21).Value = 11.5473441108545
$ws.Cells.Item(7, 23).Value = -5.464480874316946
$ws.Cells.Item(7, 24).Value = 3.623333333333334
$ws.Cells.Item(7, 25).Value = -1.98300283286118

$ws.Cells.Item(8, 1).Value = 45047
$ws.Cells.Item(8, 2).Value = 22539.418
$ws.Cells.Item(8, 3).Value = 303.365
$ws.Cells.Item(8, 4).Value = 3.7
$ws.Cells.Item(8, 5).Value = 5.06
$ws.Cells.Item(8, 6).Value = 3.57
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 9).Value = 0.6069738859241891
$ws.Cells.Item(8, 11).Value = 0.109889384619466
$ws.Cells.Item(8, 12).Value = 302.7136666666667
$ws.Cells.Item(8, 13).Value = 1.001811184061596
$ws.Cells.Item(8, 15).Value = 8.823529411764719
$ws.Cells.Item(8, 16).Value = 3.533333333333333
$ws.Cells.Item(8, 17).Value = 8.823529411764719
$ws.Cells.Item(8, 19).Value = 4.761904761904745
$ws.Cells.Item(8, 20).Value = 4.846666666666667
$ws.Cells.Item(8, 21).Value = 16.85912240184757
$ws.Cells.Item(8, 23).Value = 3.179190751445082
$ws.Cells.Item(8, 24).Value = 3.563333333333334
$ws.Cells.Item(8, 25).Value = 1.133144475920678

$ws.Cells.Item(9, 1).Value = 45078
$ws.Cells.Item(9, 2).Value = 22539.418
$ws.Cells.Item(9, 3).Value = 304.003
$ws.Cells.Item(9, 4).Value = 3.6
$ws.Cells.Item(9, 5).Value = 5.08
$ws.Cells.Item(9, 6).Value = 3.75
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 9).Value = 0.6069738859241891
$ws.Cells.Item(9, 11).Value = 0.2103077151286437
$ws.Cells.Item(9, 12).Value = 303.4666666666666
$ws.Cells.Item(9, 13).Value = 1.214225785401313
$ws.Cells.Item(9, 15).Value = -2.702702702702708
$ws.Cells.Item(9, 16).Value = 3.566666666666666
$ws.Cells.Item(9, 17).Value = 5.882352941176472
$ws.Cells.Item(9, 19).Value = 0.3952569169960674
$ws.Cells.Item(9, 20).Value = 4.989999999999999
$ws.Cells.Item(9, 21).Value = 17.32101616628174
$ws.Cells.Item(9, 23).Value = 5.042016806722693
$ws.Cells.Item(9, 24).Value = 3.593333333333334
$ws.Cells.Item(9, 25).Value = 6.232294617563738

$ws.Cells.Item(10, 1).Value = 45108
$ws.Cells.Item(10, 2).Value = 22780.933
$ws.Cells.Item(10, 3).Value = 304.628
$ws.Cells.Item(10, 4).Value = 3.5
$ws.Cells.Item(10, 5).Value = 5.12
$ws.Cells.Item(10, 6).Value = 3.9
$ws.Cells.Item(10, 7).Value = 1.071522787323076
$ws.Cells.Item(10, 9).Value = 1.685000536748049
$ws.Cells.Item(10, 11).Value = 0.2055900764137197
$ws.Cells.Item(10, 12).Value = 303.9986666666666
$ws.Cells.Item(10, 13).Value = 1.42231218953508
$ws.Cells.Item(10, 15).Value = -2.777777777777779
$ws.Cells.Item(10, 16).Value = 3.6
$ws.Cells.Item(10, 17).Value = 2.941176470588247
$ws.Cells.Item(10, 19).Value = 0.7874015748031482
$ws.Cells.Item(10, 20).Value = 5.086666666666667
$ws.Cells.Item(10, 21).Value = 18.24480369515011
$ws.Cells.Item(10, 23).Value = 4.000000000000004
$ws.Cells.Item(10, 24).Value = 3.74
$ws.Cells.Item(10, 25).Value = 10.48158640226628

$ws.Cells.Item(11, 1).Value = 45139
$ws.Cells.Item(11, 2).Value = 22780.933
$ws.Cells.Item(11, 3).Value = 306.187
$ws.Cells.Item(11, 4).Value = 3.8
$ws.Cells.Item(11, 5).Value = 5.33
$ws.Cells.Item(11, 6).Value = 4.17
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 9).Value = 1.685000536748049
$ws.Cells.Item(11, 11).Value = 0.5117717347059347
$ws.Cells.Item(11, 12).Value = 304.9393333333333
$ws.Cells.Item(11, 13).Value = 1.941362916006351
$ws.Cells.Item(11, 15).Value = 8.571428571428562
$ws.Cells.Item(11, 16).Value = 3.633333333333333
$ws.Cells.Item(11, 17).Value = 11.76470588235294
$ws.Cells.Item(11, 19).Value = 4.1015625
$ws.Cells.Item(11, 20).Value = 5.176666666666667
$ws.Cells.Item(11, 21).Value = 23.09468822170901
$ws.Cells.Item(11, 23).Value = 6.923076923076921
$ws.Cells.Item(11, 24).Value = 3.94
$ws.Cells.Item(11, 25).Value = 18.13031161473089

$ws.Cells.Item(12, 1).Value = 45170
$ws.Cells.Item(12, 2).Value = 22780.933
$ws.Cells.Item(12, 3).Value = 307.288
$ws.Cells.Item(12, 4).Value = 3.8
$ws.Cells.Item(12, 5).Value = 5.33
$ws.Cells.Item(12, 6).Value = 4.38
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 9).Value = 1.685000536748049
$ws.Cells.Item(12, 11).Value = 0.3595841756834917
$ws.Cells.Item(12, 12).Value = 306.0343333333333
$ws.Cells.Item(12, 13).Value = 2.307927925528386
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 3.7
$ws.Cells.Item(12, 17).Value = 11.76470588235294
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 5.26
$ws.Cells.Item(12, 21).Value = 23.09468822170901
$ws.Cells.Item(12, 23).Value = 5.035971223021574
$ws.Cells.Item(12, 24).Value = 4.149999999999999
$ws.Cells.Item(12, 25).Value = 24.07932011331446

$ws.Cells.Item(13, 1).Value = 45200
$ws.Cells.Item(13, 2).Value = 22960.6
$ws.Cells.Item(13, 3).Value = 307.531
$ws.Cells.Item(13, 4).Value = 3.8
$ws.Cells.Item(13, 5).Value = 5.33
$ws.Cells.Item(13, 6).Value = 4.8
$ws.Cells.Item(13, 7).Value = 0.788672702737836
$ws.Cells.Item(13, 9).Value = 2.486962378760205
$ws.Cells.Item(13, 11).Value = 0.07907890968732456
$ws.Cells.Item(13, 12).Value = 307.002
$ws.Cells.Item(13, 13).Value = 2.388831919455581
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 3.8
$ws.Cells.Item(13, 17).Value = 11.76470588235294
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(13, 20).Value = 5.33
$ws.Cells.Item(13, 21).Value = 23.09468822170901
$ws.Cells.Item(13, 23).Value = 9.589041095890405
$ws.Cells.Item(13, 24).Value = 4.449999999999999
$ws.Cells.Item(13, 25).Value = 35.97733711048159

$ws.Cells.Item(14, 1).Value = 45231
$ws.Cells.Item(14, 2).Value = 22960.6
$ws.Cells.Item(14, 3).Value = 308.024
$ws.Cells.Item(14, 4).Value = 3.7
$ws.Cells.Item(14, 5).Value = 5.33
$ws.Cells.Item(14, 6).Value = 4.5
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 9).Value = 2.486962378760205
$ws.Cells.Item(14, 11).Value = 0.160309042015272
$ws.Cells.Item(14, 12).Value = 307.6143333333333
$ws.Cells.Item(14, 13).Value = 2.552970475036287
$ws.Cells.Item(14, 15).Value = -2.631578947368407
$ws.Cells.Item(14, 16).Value = 3.766666666666667
$ws.Cells.Item(14, 17).Value = 8.823529411764719
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(14, 20).Value = 5.33
$ws.Cells.Item(14, 21).Value = 23.09468822170901
$ws.Cells.Item(14, 23).Value = -6.25
$ws.Cells.Item(14, 24).Value = 4.56
$ws.Cells.Item(14, 25).Value = 27.47875354107649

$ws.Cells.Item(15, 1).Value = 45261
$ws.Cells.Item(15, 2).Value = 22960.6
$ws.Cells.Item(15, 3).Value = 308.742
$ws.Cells.Item(15, 4).Value = 3.7
$ws.Cells.Item(15, 5).Value = 5.33
$ws.Cells.Item(15, 6).Value = 4.02
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 9).Value = 2.486962378760205
$ws.Cells.Item(15, 11).Value = 0.2330987195803047
$ws.Cells.Item(15, 12).Value = 308.099
$ws.Cells.Item(15, 13).Value = 2.792020136105156
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 3.733333333333333
$ws.Cells.Item(15, 17).Value = 8.823529411764719
$ws.Cells.Item(15, 19).Value = 0
$ws.Cells.Item(15, 20).Value = 5.33
$ws.Cells.Item(15, 21).Value = 23.09468822170901
$ws.Cells.Item(15, 23).Value = -10.66666666666668
$ws.Cells.Item(15, 24).Value = 4.44
$ws.Cells.Item(15, 25).Value = 13.88101983002832
